$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping: row -> new C value, new E value
$updates = @{
    2   = @{ C = 100826; E = 327347554 }
    63  = @{ C = 14361;  E = 36190509 }
    74  = @{ C = 950;    E = 4266554 }
    83  = @{ C = 3412;   E = 115697749 }
    92  = @{ C = 409127; E = 1595300380 }
    93  = @{ C = 209592; E = 1309153481 }
    94  = @{ C = 94195;  E = 917486940 }
    95  = @{ C = 50768;  E = 932648032 }
    96  = @{ C = 17287;  E = 793812185 }
    99  = @{ C = 165;    E = 25935853 }
    100 = @{ C = 413;    E = 66473777 }
    104 = @{ C = 135242; E = 272218580 }
    175 = @{ C = 80783;  E = 486175817 }
    177 = @{ C = 14718;  E = 251572342 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}
